$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "94.626.62"
$ws.Range("E2").Value = "  +2.82%  "
$ws.Range("D3").Value = "3.104.93"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.21"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "613.68"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.12"
$ws.Range("E7").Value = "  +2.13%  "
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.835"
$ws.Range("E10").Value = "  +14.56%  "
$ws.Range("D11").Value = "3.103.04"
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000245"
$ws.Range("E13").Value = "  -2.24%  "
$ws.Range("B14").Value = "WrappedBTC"
$ws.Range("C14").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D14").Value = "94.254.98"
$ws.Range("E14").Value = "  +2.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.60"
$ws.Range("E15").Value = "  +0.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.36"
$ws.Range("E16").Value = "  -2.54%  "
$ws.Range("D17").Value = "3.681.78"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").Value = "3.104.54"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.67"
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.85"
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.93"
$ws.Range("E21").Value = "  +1.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "450.65"
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("E23").Value = "  -1.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.98"
$ws.Range("E24").Value = "  -3.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.28"
$ws.Range("E25").Value = "  +5.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.64"
$ws.Range("E26").Value = "  +0.67%  "
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.18"
$ws.Range("E27").Value = "  +4.26%  "
$ws.Range("B28").Value = "Litecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "86.00"
$ws.Range("E28").Value = "  +6.17%  "
$ws.Range("D29").Value = "3.283.74"
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.255"
$ws.Range("E31").Value = "  +8.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.182"
$ws.Range("E32").Value = "  +9.69%  "
$ws.Range("E33").Value = "  -8.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.32"
$ws.Range("E34").Value = "  +1.54%  "
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("E36").Value = "  -1.00%  "
$ws.Range("E37").Value = "  -1.35%  "
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("E39").Value = "  -0.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.453"
$ws.Range("E40").Value = "  +4.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.99"
$ws.Range("E41").Value = "  +7.56%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.29"
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "472.78"
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.70"
$ws.Range("E44").Value = "  -12.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.26"
$ws.Range("E45").Value = "  -4.31%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "160.37"
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.691"
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("E49").Value = "  -2.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.42"
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.80"
$ws.Range("E51").Value = "  -0.55%  "
